# Refresh the scraped cryptocurrency Price (D) / Volume(1h) (E) columns
# on the active sheet for rows 2-51 with the newest run's figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.895.63"
$ws.Range("E2").Value = "  -0.02%  "

$ws.Range("D3").Value = "1.635.94"
$ws.Range("E3").Value = "  -0.66%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").Value = "'212.17"
$ws.Range("E5").Value = "  -0.64%  "

$ws.Range("D6").Value = "'0.522"
$ws.Range("E6").Value = "  -0.53%  "

$ws.Range("E7").Value = "  -0.03%  "

$ws.Range("E8").Value = "  -1.32%  "

$ws.Range("D9").Value = "'0.258"
$ws.Range("E9").Value = "  -2.89%  "

$ws.Range("E10").Value = "  -0.28%  "

$ws.Range("E11").Value = "  +0.99%  "

$ws.Range("D12").Value = "1.867.66"
$ws.Range("E12").Value = "  -0.68%  "

$ws.Range("D13").Value = "1.635.63"
$ws.Range("E13").Value = "  -0.74%  "

$ws.Range("E14").Value = "  -0.48%  "

$ws.Range("E15").Value = "  +0.89%  "

$ws.Range("D16").Value = "'65.22"
$ws.Range("E16").Value = "  -0.65%  "

$ws.Range("D17").Value = "27.891.35"
$ws.Range("E17").Value = "  -0.03%  "

$ws.Range("D18").Value = "'229.80"
$ws.Range("E18").Value = "  -0.91%  "

$ws.Range("E19").Value = "  -0.47%  "

$ws.Range("E20").Value = "  -1.99%  "

$ws.Range("E21").Value = "  -0.11%  "

$ws.Range("E22").Value = "  -0.40%  "

$ws.Range("D23").Value = "'10.33"
$ws.Range("E23").Value = "  -3.55%  "

$ws.Range("E24").Value = "  -3.27%  "

$ws.Range("D25").Value = "'153.18"
$ws.Range("E25").Value = "  +0.72%  "

$ws.Range("E26").Value = "  +0.65%  "

$ws.Range("E27").Value = "  -0.63%  "

$ws.Range("E28").Value = "  -0.61%  "

$ws.Range("E29").Value = "  +0.00%  "

$ws.Range("E30").Value = "  -0.97%  "

$ws.Range("E31").Value = "  -0.79%  "

$ws.Range("E32").Value = "  +0.54%  "

$ws.Range("D33").Value = "1.408.15"
$ws.Range("E33").Value = "  -3.21%  "

$ws.Range("E34").Value = "  -1.52%  "

$ws.Range("D35").Value = "'1.57"
$ws.Range("E35").Value = "  +1.06%  "

$ws.Range("D36").Value = "'1.00"
$ws.Range("E36").Value = "  +8.70%  "

$ws.Range("E37").Value = "  +1.53%  "

$ws.Range("E38").Value = "  +0.41%  "

$ws.Range("D39").Value = "'0.562"
$ws.Range("E39").Value = "  -0.03%  "

$ws.Range("D40").Value = "'0.872"
$ws.Range("E40").Value = "  -1.93%  "

$ws.Range("E41").Value = "  +0.28%  "

$ws.Range("E42").Value = "  -0.09%  "

$ws.Range("D43").Value = "'66.92"
$ws.Range("E43").Value = "  -3.48%  "

$ws.Range("D44").Value = "'5.50"
$ws.Range("E44").Value = "  +2.35%  "

$ws.Range("E45").Value = "  +1.35%  "

$ws.Range("E46").Value = "  -1.60%  "

$ws.Range("D47").Value = "1.777.21"
$ws.Range("E47").Value = "  -0.73%  "

$ws.Range("D48").Value = "'87.76"
$ws.Range("E48").Value = "  -1.36%  "

$ws.Range("E49").Value = "  -0.49%  "

$ws.Range("E50").Value = "  -0.20%  "

$ws.Range("D51").Value = "'7.55"
$ws.Range("E51").Value = "  -2.68%  "
